$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new "Pub Date" column before the current "Media Type(s)" column (G) ---
$ws.Columns("G").Insert()

$ws.Range("G1").Value = "Pub Date"
$ws.Range("G2").Value = 2025
$ws.Range("G4").Value = 2025
$ws.Range("G7").Value = 2025

# --- 2. Update the "Media Type(s)" text for the two JPEG Trust rows (now column H) ---
$ws.Range("H4").Value = "Any (image focused)"
$ws.Range("H5").Value = "Any (image focused)"

# --- 3. Add the new "TDM Reservation Protocol" row (row 8) ---
# Copy formatting from an existing similarly-formatted row first so the new
# row reuses existing cell styles instead of minting new ones.
$ws.Range("A6:F6").Copy($ws.Range("A8:F8"))
$ws.Range("A6").Copy($ws.Range("G8:I8"))

$ws.Range("A8").Value = "TDM Reservation Protocol"
$ws.Range("B8").Value = "Rights Declarations"
$ws.Range("C8").Value = "W3C"
$ws.Range("D8").Value = "TDMRep"
$ws.Range("E8").Value = "https://www.w3.org/ns/tdmrep/"
$ws.Range("F8").Value = "Published"
$ws.Range("G8").ClearContents() | Out-Null
$ws.Range("H8").Value = "Web pages`nEPUB`nPDF"
$ws.Range("I8").Value = "This protocol provides guidelines for reserving content from text and data mining. It includes methods for creating and maintaining TDMRep files, which can be used to document the reservation of digital assets. This helps in ensuring that content is not used for data mining without the creator's consent."

# --- 4. Row heights ---
$ws.Rows(1).RowHeight = 40
$ws.Rows(2).RowHeight = 60
$ws.Rows(3).RowHeight = 40
$ws.Rows(4).RowHeight = 100
$ws.Rows(5).RowHeight = 40
$ws.Rows(6).RowHeight = 40
$ws.Rows(7).RowHeight = 380
$ws.Rows(8).RowHeight = 409.6

# --- 5. Column width for the new Pub Date column ---
$ws.Range("G1").ColumnWidth = $ws.Range("F1").ColumnWidth

# --- 6. Selection / view ---
$ws.Range("H6").Select() | Out-Null
